$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" text note in A1 with the new conversion rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.79 = 6616.3 pesos`n✅ 6616.3 pesos = 1.78 = 935.31 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update the numeric rate values in the "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 560
$wsTasas.Range("O10").Value = 3705.13
$wsTasas.Range("N12").Value = 3713.8
$wsTasas.Range("O12").Value = 525
